# Append newly scraped Lancers listings and refresh the "fetched at" timestamps.
# New timestamp for this scrape run:
#   2025-10-26 12:33:22
#
# The three brand-new postings (顧客予約サイン / 楽天配布クーポン / IB報酬EA) are
# inserted right after the existing top row, pushing the previously-second row
# ("Power Automate for Desktop...") down to row 6, followed by two more new
# postings (進行管理 / eBayテラピーク).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-10-26 12:33:22"

# --- Row 2: existing top item - only the fetch timestamp changes ---
$ws.Range("A2").Value = $newTimestamp

# --- Rows 3-8: final desired content for this block of the sheet ---
$rows = @(
    @{ Row = 3;  A = $newTimestamp; B = "【システム開発】顧客予約サインシステムの構築依頼"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定";   E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5420832"; G = 113; H = "◆開発,システム開発" },
    @{ Row = 4;  A = $newTimestamp; B = "【急募】楽天配布型クーポン登録の自動化システム構築";                 C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定";   E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5420867"; G = 98;  H = "◆自動化" },
    @{ Row = 5;  A = $newTimestamp; B = "IB報酬を得るための高性能EA開発依頼";                                 C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5420753"; G = 68;  H = "◆開発" },
    @{ Row = 6;  A = $newTimestamp; B = "【Power Automate for Desktop】販売管理システムへExcelから自動入力"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定";   E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5407216"; G = 48;  H = "◇管理" },
    @{ Row = 7;  A = $newTimestamp; B = "進行管理およびチームディレクションを担当";                           C = "システム開発"; D = "~ 5,000 円 / 固定";               E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5418064"; G = 30;  H = "◇管理" },
    @{ Row = 8;  A = $newTimestamp; B = "eBayテラピークでのキーワード検索結果等の取得するためのシステム制作"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定";   E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5420779"; G = 33;  H = $null }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    if ($item.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $item.H
    }
}

# --- Hyperlinks: rebuild F2:F8 cleanly and in order so relationship ids line up ---
# (Range.Hyperlinks.Delete() clears every hyperlink on the sheet in this engine,
# so just clear once and re-add all seven in top-to-bottom order.)
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5420678") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5420832") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5420867") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5420753") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5407216") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5418064") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5420779") | Out-Null

# Re-apply the standard "Hyperlink" cell style on each link cell (Hyperlinks.Add
# leaves behind an extra style record; forcing the named style keeps every
# linked cell on the workbook's single canonical Hyperlink style).
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("F8").Style = "Hyperlink"

Write-Output "Appended scrape batch for $newTimestamp (rows 3-8); refreshed row 2 timestamp."
